$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" '64.730.91'
Set-TextValue "E2" '  +2.01%  '
Set-TextValue "D3" '2.634.19'
Set-TextValue "E3" '  +2.14%  '
Set-TextValue "E4" '  -0.02%  '
Set-TextValue "D5" '592.69'
Set-TextValue "E5" '  +0.63%  '
Set-TextValue "D6" '154.78'
Set-TextValue "E6" '  +3.10%  '
Set-TextValue "E7" '  -0.01%  '
Set-TextValue "D8" '0.589'
Set-TextValue "E8" '  +0.43%  '
Set-TextValue "D9" '0.116'
Set-TextValue "E9" '  +5.17%  '
Set-TextValue "D10" '0.395'
Set-TextValue "E10" '  +3.16%  '
Set-TextValue "E11" '  +0.68%  '
Set-TextValue "E12" '  +1.78%  '
Set-TextValue "D13" '28.82'
Set-TextValue "E13" '  +4.65%  '
Set-TextValue "E14" '  +18.54%  '
Set-TextValue "D15" '3.103.24'
Set-TextValue "E15" '  +1.93%  '
Set-TextValue "D16" '64.683.86'
Set-TextValue "E16" '  +2.17%  '
Set-TextValue "D17" '2.625.16'
Set-TextValue "E17" '  +2.42%  '
Set-TextValue "D18" '12.50'
Set-TextValue "E18" '  +2.58%  '
Set-TextValue "D19" '4.76'
Set-TextValue "E19" '  +0.94%  '
Set-TextValue "D20" '349.69'
Set-TextValue "E20" '  +0.91%  '
Set-TextValue "D21" '7.25'
Set-TextValue "E21" '  +5.79%  '
Set-TextValue "E22" '  -0.10%  '
Set-TextValue "D23" '67.81'
Set-TextValue "E23" '  +0.84%  '
Set-TextValue "D24" '1.68'
Set-TextValue "E24" '  -0.67%  '
Set-TextValue "D25" '9.45'
Set-TextValue "E25" '  +3.30%  '
Set-TextValue "B27" 'Aptos'
Set-TextValue "C27" 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D27" '8.08'
Set-TextValue "E27" '  +0.06%  '
Set-TextValue "B28" 'Kaspa'
Set-TextValue "C28" 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D28" '0.163'
Set-TextValue "E28" '  +0.12%  '
Set-TextValue "B29" 'Binance-PegBSC-USD'
Set-TextValue "C29" 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
Set-TextValue "D29" '0.993'
Set-TextValue "E29" '  -0.63%  '
Set-TextValue "D30" '0.0₃0927'
Set-TextValue "E30" '  +7.65%  '
Set-TextValue "D31" '2.07'
Set-TextValue "E31" '  +1.63%  '
Set-TextValue "D32" '509.31'
Set-TextValue "E32" '  -8.15%  '
Set-TextValue "E33" '  +0.06%  '
Set-TextValue "E34" '  +6.83%  '
Set-TextValue "D35" '6.20'
Set-TextValue "E35" '  +2.17%  '
Set-TextValue "E36" '  +2.23%  '
Set-TextValue "D37" '164.73'
Set-TextValue "E37" '  -1.15%  '
Set-TextValue "D38" '20.03'
Set-TextValue "E38" '  +2.45%  '
Set-TextValue "E39" '  +4.04%  '
Set-TextValue "D40" '0.999'
Set-TextValue "E40" '  -0.05%  '
Set-TextValue "D41" '1.00'
Set-TextValue "E41" '  +0.07%  '
Set-TextValue "D42" '42.17'
Set-TextValue "E42" '  +6.20%  '
Set-TextValue "D43" '163.50'
Set-TextValue "E43" '  -1.36%  '
Set-TextValue "E44" '  +2.01%  '
Set-TextValue "D45" '0.0608'
Set-TextValue "E45" '  +3.67%  '
Set-TextValue "D46" '22.67'
Set-TextValue "E46" '  -1.24%  '
Set-TextValue "E47" '  +3.72%  '
Set-TextValue "D48" '0.643'
Set-TextValue "E48" '  +2.41%  '
Set-TextValue "E49" '  +0.60%  '
Set-TextValue "E50" '  +1.50%  '
Set-TextValue "D51" '19.21'
Set-TextValue "E51" '  +0.53%  '
